# Auto-generated Excel COM-interop script applying the Masamune_Profits value updates.
# Each row entry sets the changed numeric cells and clears any cells that were
# removed entirely in the target revision.

$wb = $excel.ActiveWorkbook

function Set-Cells($ws, $row, $setMap, $clearCols) {
    foreach ($col in $setMap.Keys) {
        $ws.Range("$col$row").Value = $setMap[$col]
    }
    foreach ($col in $clearCols) {
        $ws.Range("$col$row").ClearContents()
    }
}

$ws = $wb.Worksheets.Item("ALC")
Set-Cells $ws 101 @{"H" = 391.375; "I" = 411.7143; "J" = 249; "K" = 1235.1429; "L" = 747; "M" = 386.8571000000002; "N" = -3991} @()
Set-Cells $ws 103 @{"H" = 1097.3572; "I" = 859.1429000000001; "J" = 1335.5714; "K" = 2577.4287; "L" = 4006.7142; "M" = -1991.4287; "N" = -5178.7142} @()
Set-Cells $ws 120 @{"H" = 49718; "J" = 49718; "L" = 49718; "N" = -59394} @()
Set-Cells $ws 129 @{"H" = 1517.4814; "I" = 2211; "K" = 6633; "M" = -1633} @()
Set-Cells $ws 130 @{"H" = 49772; "J" = 49772; "L" = 49772; "N" = -59812} @()
Set-Cells $ws 132 @{"H" = 13546.595; "I" = 2193.0923; "J" = 66259.28999999999; "K" = 6579.276899999999; "L" = 198777.87; "M" = -4049.276899999999; "N" = -203837.87} @()
$ws = $wb.Worksheets.Item("ARM")
Set-Cells $ws 32 @{"H" = 8502.8125; "I" = 7988.6377; "J" = 11728.091; "K" = 7988.6377; "L" = 11728.091; "M" = -7701.6377; "N" = -12302.091} @()
Set-Cells $ws 45 @{"H" = 1781.32; "I" = 1717.625; "J" = 1894.5555; "K" = 1717.625; "L" = 1894.5555; "M" = -1340.625; "N" = -2648.5555} @()
Set-Cells $ws 74 @{"H" = 1345.5352; "I" = 1383.875; "J" = 1265.5217; "K" = 1383.875; "L" = 1265.5217; "M" = -509.875; "N" = -3013.5217} @()
Set-Cells $ws 77 @{"H" = 1345.5352; "I" = 1383.875; "J" = 1265.5217; "K" = 6919.375; "L" = 6327.6085; "M" = -2551.375; "N" = -15063.6085} @()
Set-Cells $ws 80 @{"H" = 45632.668; "J" = 45632.668; "L" = 45632.668; "N" = -47628.668} @()
Set-Cells $ws 83 @{"H" = 45632.668; "J" = 45632.668; "L" = 136898.004; "N" = -146882.004} @()
Set-Cells $ws 111 @{"H" = 49543; "J" = 49543; "L" = 49543; "N" = -57723} @()
Set-Cells $ws 122 @{"H" = 2247.2942; "I" = 2361.4666; "J" = 1391; "K" = 7084.399800000001; "L" = 4173; "M" = -4634.399800000001; "N" = -9073} @()
Set-Cells $ws 131 @{"H" = 47642.8; "J" = 47642.8; "L" = 47642.8; "N" = -57722.8} @()
Set-Cells $ws 132 @{"H" = 10871268; "I" = 16130152; "J" = 2907.6; "K" = 48390456; "L" = 8722.799999999999; "M" = -48387926; "N" = -13782.8} @()
$ws = $wb.Worksheets.Item("BSM")
Set-Cells $ws 20 @{"H" = 2009.75; "I" = 1713.3125; "J" = 2405; "K" = 1713.3125; "L" = 2405; "M" = -1466.3125; "N" = -2899} @()
Set-Cells $ws 92 @{"H" = 39501.8; "J" = 39501.8; "L" = 39501.8; "N" = -44493.8} @()
Set-Cells $ws 105 @{"H" = 1832.8928; "I" = 1568.3334; "J" = 2138.1538; "K" = 1568.3334; "L" = 2138.1538; "M" = 178.6666; "N" = -5632.1538} @()
Set-Cells $ws 110 @{"H" = 46400; "J" = 46400; "L" = 46400; "N" = -54580} @()
Set-Cells $ws 124 @{"H" = 43998.2; "J" = 43998.2; "L" = 43998.2; "N" = -53818.2} @()
Set-Cells $ws 134 @{"H" = 2522.8494; "I" = 1621.4651; "J" = 3298.04; "K" = 4864.3953; "L" = 9894.119999999999; "M" = -2329.3953; "N" = -14964.12} @()
$ws = $wb.Worksheets.Item("CRP")
Set-Cells $ws 31 @{"H" = 2662.17; "I" = 1123.7812; "J" = 3386.1177; "K" = 1123.7812; "L" = 3386.1177; "M" = -828.7811999999999; "N" = -3976.1177} @()
Set-Cells $ws 34 @{"H" = 2662.17; "I" = 1123.7812; "J" = 3386.1177; "K" = 1123.7812; "L" = 3386.1177; "M" = -921.7811999999999; "N" = -3790.1177} @()
Set-Cells $ws 107 @{"H" = 564; "I" = 513.2; "K" = 513.2; "M" = 1406.8} @()
Set-Cells $ws 134 @{"H" = 468322.12; "I" = 1206.7; "J" = 1402553; "K" = 3620.1; "L" = 4207659; "M" = -1085.1; "N" = -4212729} @()
$ws = $wb.Worksheets.Item("CUL")
Set-Cells $ws 113 @{"H" = 3279.3513; "I" = 4926.2173; "J" = 573.7857; "K" = 14778.6519; "L" = 1721.3571; "M" = -12608.6519; "N" = -6061.3571} @()
Set-Cells $ws 131 @{"H" = 9686.846; "I" = 20687.8; "J" = 2811.25; "K" = 62063.39999999999; "L" = 8433.75; "M" = -57023.39999999999; "N" = -18513.75} @()
Set-Cells $ws 136 @{"H" = 29415164; "I" = 55557696; "J" = 4812.5; "K" = 166673088; "L" = 14437.5; "M" = -166667988; "N" = -24637.5} @()
$ws = $wb.Worksheets.Item("GSM")
Set-Cells $ws 92 @{"H" = 4808.1; "I" = 4250; "J" = 4947.625; "K" = 4250; "L" = 4947.625; "M" = -2378; "N" = -8691.625} @()
Set-Cells $ws 113 @{"H" = 1430.9333; "I" = 1633.1666; "K" = 1633.1666; "M" = 536.8334} @()
Set-Cells $ws 122 @{"H" = 1390.7; "I" = 1378.5555; "J" = 1500; "K" = 4135.666499999999; "L" = 4500; "M" = -1685.666499999999; "N" = -9400} @()
Set-Cells $ws 126 @{"H" = 17687.715; "I" = 23382.8; "K" = 70148.39999999999; "M" = -67678.39999999999} @()
Set-Cells $ws 130 @{"H" = 37957; "J" = 37957; "L" = 37957; "N" = -47997} @()
Set-Cells $ws 132 @{"H" = 1952.4; "I" = 923.6; "J" = 3667.0667; "K" = 2770.8; "L" = 11001.2001; "M" = -240.8000000000002; "N" = -16061.2001} @()
$ws = $wb.Worksheets.Item("LTW")
Set-Cells $ws 7 @{"H" = 2287.9268; "I" = 1788.2354; "J" = 4715; "K" = 1788.2354; "L" = 4715; "M" = -1676.2354; "N" = -4939} @()
Set-Cells $ws 31 @{"H" = 568.75; "I" = 568.75; "J" = 0; "K" = 568.75; "L" = 0; "M" = -320.75} @("N")
Set-Cells $ws 46 @{"H" = 9272.666999999999; "I" = 2350; "J" = 10657.2; "K" = 2350; "L" = 10657.2; "M" = -2162; "N" = -11033.2} @()
Set-Cells $ws 68 @{"H" = 3500; "I" = 0; "J" = 3500; "K" = 0; "L" = 3500; "N" = -4998} @("M")
Set-Cells $ws 71 @{"H" = 3500; "I" = 0; "J" = 3500; "K" = 0; "L" = 17500; "N" = -24988} @("M")
Set-Cells $ws 100 @{"H" = 1379.625; "I" = 1404.9333; "J" = 1000; "K" = 1404.9333; "L" = 1000; "M" = -863.9332999999999; "N" = -2082} @()
Set-Cells $ws 126 @{"H" = 2287.9268; "I" = 1788.2354; "J" = 4715; "K" = 5364.706200000001; "L" = 14145; "M" = -2894.706200000001; "N" = -19085} @()
$ws = $wb.Worksheets.Item("WVR")
Set-Cells $ws 32 @{"H" = 24205.8; "I" = 15000; "J" = 26507.25; "K" = 15000; "L" = 26507.25; "M" = -14683; "N" = -27141.25} @()
Set-Cells $ws 62 @{"H" = 3249.1904; "I" = 3333; "K" = 3333; "M" = -2709} @()
Set-Cells $ws 65 @{"H" = 3249.1904; "I" = 3333; "K" = 16665; "M" = -13545} @()
Set-Cells $ws 92 @{"H" = 28023; "J" = 28023; "L" = 28023; "N" = -33015} @()
Set-Cells $ws 132 @{"H" = 1583.8795; "I" = 1299.7188; "J" = 2541.0527; "K" = 3899.1564; "L" = 7623.158100000001; "M" = -1369.1564; "N" = -12683.1581} @()
